$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.315959215164185
$ws.Range("B1").Value = -1
$ws.Range("D1").Value = 1.34495210647583
$ws.Range("E1").Value = 0.8169378638267517
